$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 894.2
$ws.Range("I15").Value = 894.2
$ws.Range("K15").Value = 2682.6
$ws.Range("M15").Value = -2513.6
$ws.Range("H17").Value = 1084177.5
$ws.Range("J17").Value = 1113084
$ws.Range("L17").Value = 3339252
$ws.Range("N17").Value = -3339588
$ws.Range("H29").Value = 70000
$ws.Range("J29").Value = 68333.336
$ws.Range("L29").Value = 205000.008
$ws.Range("N29").Value = -205562.008
$ws.Range("H32").Value = 1767.6
$ws.Range("I32").Value = 1645
$ws.Range("J32").Value = 1798.25
$ws.Range("K32").Value = 1645
$ws.Range("L32").Value = 1798.25
$ws.Range("M32").Value = -1319
$ws.Range("N32").Value = -2450.25
$ws.Range("H38").Value = 1294.2354
$ws.Range("I38").Value = 425.25
$ws.Range("J38").Value = 2066.6667
$ws.Range("K38").Value = 1275.75
$ws.Range("L38").Value = 6200.000100000001
$ws.Range("M38").Value = -903.75
$ws.Range("N38").Value = -6944.000100000001
$ws.Range("H43").Value = 10987.363
$ws.Range("I43").Value = 5486.5
$ws.Range("J43").Value = 12209.777
$ws.Range("K43").Value = 5486.5
$ws.Range("L43").Value = 12209.777
$ws.Range("M43").Value = -5417.5
$ws.Range("N43").Value = -12347.777
$ws.Range("H51").Value = 3598.182
$ws.Range("I51").Value = 2793.3333
$ws.Range("K51").Value = 2793.3333
$ws.Range("M51").Value = -2309.3333
$ws.Range("H58").Value = 969.7143
$ws.Range("J58").Value = 1998
$ws.Range("L58").Value = 5994
$ws.Range("N58").Value = -6294
$ws.Range("H111").Value = 2099.2222
$ws.Range("J111").Value = 1166
$ws.Range("L111").Value = 3498
$ws.Range("N111").Value = -9632
$ws.Range("H132").Value = 2553.4443
$ws.Range("I132").Value = 2594.3901
$ws.Range("J132").Value = 2424.3076
$ws.Range("K132").Value = 7783.1703
$ws.Range("L132").Value = 7272.9228
$ws.Range("M132").Value = -5253.1703
$ws.Range("N132").Value = -12332.9228
$ws.Range("H134").Value = 37744.844
$ws.Range("J134").Value = 40899.63
$ws.Range("L134").Value = 40899.63
$ws.Range("N134").Value = -51039.63
$ws.Range("H135").Value = 1560.4043
$ws.Range("I135").Value = 1111.1945
$ws.Range("J135").Value = 3030.5454
$ws.Range("K135").Value = 10000.7505
$ws.Range("L135").Value = 27274.9086
$ws.Range("M135").Value = -7465.7505
$ws.Range("N135").Value = -32344.9086
$ws.Range("H138").Value = 2130.6826
$ws.Range("I138").Value = 1212.1818
$ws.Range("J138").Value = 3141.0334
$ws.Range("K138").Value = 3636.5454
$ws.Range("L138").Value = 9423.100199999999
$ws.Range("M138").Value = 1503.4546
$ws.Range("N138").Value = -19703.1002
$ws.Range("H139").Value = 39355.6
$ws.Range("J139").Value = 39355.6
$ws.Range("L139").Value = 39355.6
$ws.Range("N139").Value = -49635.6
$ws.Range("H141").Value = 1870.2
$ws.Range("I141").Value = 1278.8096
$ws.Range("K141").Value = 3836.4288
$ws.Range("M141").Value = 1343.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8203.736999999999
$ws.Range("I32").Value = 7535.5737
$ws.Range("K32").Value = 7535.5737
$ws.Range("M32").Value = -7248.5737
$ws.Range("H61").Value = 20834906
$ws.Range("I61").Value = 25001658
$ws.Range("J61").Value = 1145
$ws.Range("K61").Value = 25001658
$ws.Range("L61").Value = 1145
$ws.Range("M61").Value = -25001446
$ws.Range("N61").Value = -1569
$ws.Range("H74").Value = 8930862
$ws.Range("I74").Value = 13159277
$ws.Range("J74").Value = 4205.778
$ws.Range("K74").Value = 13159277
$ws.Range("L74").Value = 4205.778
$ws.Range("M74").Value = -13158403
$ws.Range("N74").Value = -5953.778
$ws.Range("H77").Value = 8930862
$ws.Range("I77").Value = 13159277
$ws.Range("J77").Value = 4205.778
$ws.Range("K77").Value = 65796385
$ws.Range("L77").Value = 21028.89
$ws.Range("M77").Value = -65792017
$ws.Range("N77").Value = -29764.89
$ws.Range("H122").Value = 4682.485
$ws.Range("I122").Value = 6429.85
$ws.Range("J122").Value = 1994.2307
$ws.Range("K122").Value = 19289.55
$ws.Range("L122").Value = 5982.6921
$ws.Range("M122").Value = -16839.55
$ws.Range("N122").Value = -10882.6921
$ws.Range("H132").Value = 5683863.5
$ws.Range("I132").Value = 10871461
$ws.Range("J132").Value = 2208.476
$ws.Range("K132").Value = 32614383
$ws.Range("L132").Value = 6625.428
$ws.Range("M132").Value = -32611853
$ws.Range("N132").Value = -11685.428
$ws.Range("H136").Value = 20834906
$ws.Range("I136").Value = 25001658
$ws.Range("J136").Value = 1145
$ws.Range("K136").Value = 75004974
$ws.Range("L136").Value = 3435
$ws.Range("M136").Value = -75002424
$ws.Range("N136").Value = -8535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5745.25
$ws.Range("I134").Value = 3892.5
$ws.Range("J134").Value = 8833.166999999999
$ws.Range("K134").Value = 11677.5
$ws.Range("L134").Value = 26499.501
$ws.Range("M134").Value = -9142.5
$ws.Range("N134").Value = -31569.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7250388
$ws.Range("I31").Value = 4346.39
$ws.Range("J31").Value = 66667930
$ws.Range("K31").Value = 4346.39
$ws.Range("L31").Value = 66667930
$ws.Range("M31").Value = -4051.39
$ws.Range("N31").Value = -66668520
$ws.Range("H34").Value = 7250388
$ws.Range("I34").Value = 4346.39
$ws.Range("J34").Value = 66667930
$ws.Range("K34").Value = 4346.39
$ws.Range("L34").Value = 66667930
$ws.Range("M34").Value = -4144.39
$ws.Range("N34").Value = -66668334
$ws.Range("H122").Value = 1753.4783
$ws.Range("I122").Value = 1833.6842
$ws.Range("J122").Value = 1372.5
$ws.Range("K122").Value = 5501.0526
$ws.Range("L122").Value = 4117.5
$ws.Range("M122").Value = -3051.0526
$ws.Range("N122").Value = -9017.5
$ws.Range("H132").Value = 11365788
$ws.Range("I132").Value = 15626722
$ws.Range("J132").Value = 3295.8333
$ws.Range("K132").Value = 46880166
$ws.Range("L132").Value = 9887.499899999999
$ws.Range("M132").Value = -46877636
$ws.Range("N132").Value = -14947.4999
$ws.Range("H134").Value = 1292.3889
$ws.Range("I134").Value = 1386.4445
$ws.Range("J134").Value = 1010.2222
$ws.Range("K134").Value = 4159.333500000001
$ws.Range("L134").Value = 3030.6666
$ws.Range("M134").Value = -1624.333500000001
$ws.Range("N134").Value = -8100.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 300
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("H117").Value = 1500
$ws.Range("J117").Value = 1500
$ws.Range("L117").Value = 4500
$ws.Range("N117").Value = -11384
$ws.Range("H118").Value = 1312.9
$ws.Range("J118").Value = 1239.8948
$ws.Range("L118").Value = 3719.6844
$ws.Range("N118").Value = -6205.6844
$ws.Range("H129").Value = 3278.5
$ws.Range("J129").Value = 2736.7917
$ws.Range("L129").Value = 8210.375100000001
$ws.Range("N129").Value = -18210.3751
$ws.Range("H131").Value = 833.42267
$ws.Range("J131").Value = 859.3555
$ws.Range("L131").Value = 2578.0665
$ws.Range("N131").Value = -12658.0665
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14165.4
$ws.Range("I70").Value = 52827
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 52827
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -52557
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 14165.4
$ws.Range("I73").Value = 52827
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 52827
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -51891
$ws.Range("N73").Value = -6372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 10642491
$ws.Range("I136").Value = 15626634
$ws.Range("J136").Value = 9654
$ws.Range("K136").Value = 46879902
$ws.Range("L136").Value = 28962
$ws.Range("M136").Value = -46877352
$ws.Range("N136").Value = -34062
$ws.Range("H139").Value = 59834.855
$ws.Range("J139").Value = 59834.855
$ws.Range("L139").Value = 59834.855
$ws.Range("N139").Value = -70114.85500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2276.3462
$ws.Range("I96").Value = 1127.0714
$ws.Range("J96").Value = 3617.1667
$ws.Range("K96").Value = 1127.0714
$ws.Range("L96").Value = 3617.1667
$ws.Range("M96").Value = 245.9286
$ws.Range("N96").Value = -6363.1667
$ws.Range("H107").Value = 1290.3636
$ws.Range("I107").Value = 1786.2858
$ws.Range("J107").Value = 422.5
$ws.Range("K107").Value = 5358.857400000001
$ws.Range("L107").Value = 1267.5
$ws.Range("M107").Value = -3438.857400000001
$ws.Range("N107").Value = -5107.5
$ws.Range("H122").Value = 2049.5454
$ws.Range("I122").Value = 2049.5454
$ws.Range("K122").Value = 6148.6362
$ws.Range("M122").Value = -3698.6362
$ws.Range("H132").Value = 1769.862
$ws.Range("I132").Value = 1223.2
$ws.Range("J132").Value = 2355.5715
$ws.Range("K132").Value = 3669.6
$ws.Range("L132").Value = 7066.7145
$ws.Range("M132").Value = -1139.6
$ws.Range("N132").Value = -12126.7145
$ws.Range("H136").Value = 914.549
$ws.Range("I136").Value = 775.51166
$ws.Range("J136").Value = 1661.875
$ws.Range("K136").Value = 2326.53498
$ws.Range("L136").Value = 4985.625
$ws.Range("M136").Value = 223.4650200000001
$ws.Range("N136").Value = -10085.625
$ws.Range("H138").Value = 54549.5
$ws.Range("J138").Value = 54549.5
$ws.Range("L138").Value = 54549.5
$ws.Range("N138").Value = -64829.5
